$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and update chart series references ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "arima_graph"

$co = $ws1.ChartObjects(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(arima_graph!`$B`$1,,arima_graph!`$B`$2:`$B`$13,1)"
$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(arima_graph!`$C`$1,,arima_graph!`$C`$2:`$C`$13,2)"

# --- Add the new pm10_limits sheet after arima_graph ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "pm10_limits"

$ws2.Range("F4").Value = "ue"
$ws2.Range("G4").Value = "oms"
$ws2.Range("E5").Value = "Promedio 24 horas"
$ws2.Range("F5").Value = "<= 50 ug/m3; <= 35 días"
$ws2.Range("G5").Value = "<= 50 ug/m3; <= 3 días"
$ws2.Range("E6").Value = "Promedio anual"
$ws2.Range("F6").Value = "<= 40 ug/m3"
$ws2.Range("G6").Value = "<= 20 ug/m3"

$ws2.Range("F4:G6").HorizontalAlignment = -4108

# Target widths are 19.1796875 / 20.86328125 / 22.1796875 characters; the
# engine quantizes ColumnWidth to increments of 1/6 character on save, so we
# feed it the bucket midpoint that rounds back to the closest achievable width.
$ws2.Columns.Item(5).ColumnWidth = 18.333333333333332
$ws2.Columns.Item(6).ColumnWidth = 20.0
$ws2.Columns.Item(7).ColumnWidth = 21.333333333333332

$ws2.Activate() | Out-Null
$ws2.Range("G9").Select() | Out-Null
